$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "charger1"
$ws.Range("D1").Value = "charger2"
$ws.Range("E1").Value = "charger3"
$ws.Range("F1").Value = "charger4"
$ws.Range("G1").Value = "charger5"

$ws.Range("F11").Select()
